# Published state of ETDataset on 1 June 2018
#
# The "Dashboard" sheet (Sheet2) had five manual-input rows
# (simult_sd, simult_se, simult_wd, simult_we, peak_load_units_present)
# removed entirely - a plain full-row delete of rows 21:25, which
# shifts every row below it up by five.
$wb = $excel.ActiveWorkbook

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Rows("21:25").Select()
$dashboard.Rows("21:25").Delete()

# The "Research data" sheet became the active/selected tab (it was
# "Dashboard" before).
$wb.Worksheets.Item("Research data").Activate()
